# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" summary text ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsHoja1.Range("A1")
$text = $cellA1.Value2
$text = $text.Replace("1000 Bs = 2.01 = 7512.05 pesos", "1000 Bs = 2.04 = 7721.08 pesos")
$text = $text.Replace("7512.05 pesos = 2.0 = 925.02 Bs", "7721.08 pesos = 2.07 = 929.26 Bs")
$cellA1.Value2 = $text

# --- tasas: update the raw rate inputs (N10, O10, N12, O12) ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 489
$wsTasas.Range("O10").Value = 3775.61
$wsTasas.Range("N12").Value = 3739
$wsTasas.Range("O12").Value = 450
